# Update gene region coordinates per "Regions Data Tables" mail.
#
# 1) BMPR1A-relevant deletion region (row 3, chr10): narrow the hg38 window
#    from 121 kb to 16 kb -> hg38 chr10:86754489-86770921.
# 2) GREM1-associated CNV region (row 4, chr15): narrow the hg38 window from
#    ~400 kb to the 40 kb region from PMID: 22561515 -> hg38
#    chr15:32672738-32712558, and update the comment to cite that PMID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - BMPR1A promoter: hg38_start / hg38_stop (columns I/J)
$ws.Range("I3").Value2 = 86754489
$ws.Range("J3").Value2 = 86770921

# Row 4 - GREM1 promoter: hg38_start / hg38_stop (columns I/J) + comment (K)
$ws.Range("I4").Value2 = 32672738
$ws.Range("J4").Value2 = 32712558
$ws.Range("K4").Value2 = "GREM1 promoter (based on PMID: 22561515)"

# Leave the active selection where the edit was made, matching the saved file.
$ws.Range("J3").Select()
